$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.11250033429948
$ws.Range("C2").Value = 7.686099520365063
$ws.Range("D2").Value = 9.91673279530918
$ws.Range("F2").Value = 51.96771528803502
$ws.Range("G2").Value = 3.721937179831112
$ws.Range("L2").Value = 10.50928741664509
$ws.Range("M2").Value = 17.60381300096119
$ws.Range("B3").Value = 19.89029251911209
$ws.Range("C3").Value = 7.261925906960638
$ws.Range("D3").Value = 9.79250606301763
$ws.Range("F3").Value = 50.81244701226309
$ws.Range("G3").Value = 3.726795783208224
$ws.Range("L3").Value = 10.51759085995673
$ws.Range("M3").Value = 17.5966123249732
$ws.Range("B4").Value = 19.76188042726407
$ws.Range("C4").Value = 6.987130205838539
$ws.Range("D4").Value = 9.715481633371384
$ws.Range("F4").Value = 50.09666155632586
$ws.Range("G4").Value = 3.729927666195253
$ws.Range("L4").Value = 10.52400068067035
$ws.Range("M4").Value = 17.59715256206053
$ws.Range("B5").Value = 19.7116376339409
$ws.Range("C5").Value = 6.871564464818217
$ws.Range("D5").Value = 9.683920851796305
$ws.Range("F5").Value = 49.80368241409713
$ws.Range("G5").Value = 3.731241490870063
$ws.Range("L5").Value = 10.52694238790338
$ws.Range("M5").Value = 17.598620047257
$ws.Range("B6").Value = 19.70342274772607
$ws.Range("C6").Value = 6.85215924131874
$ws.Range("D6").Value = 9.678670231287166
$ws.Range("F6").Value = 49.75496533927074
$ws.Range("G6").Value = 3.731461923778028
$ws.Range("L6").Value = 10.52745076178379
$ws.Range("M6").Value = 17.59893903274946
$ws.Range("B7").Value = 19.76119430213284
$ws.Range("C7").Value = 6.985586120778284
$ws.Range("D7").Value = 9.715056672375741
$ws.Range("F7").Value = 50.092715119628
$ws.Range("G7").Value = 3.729945232601296
$ws.Range("L7").Value = 10.52403901900063
$ws.Range("M7").Value = 17.59716730364183
$ws.Range("B8").Value = 20.03426025848445
$ws.Range("C8").Value = 7.542835618441038
$ws.Range("D8").Value = 9.874064151512286
$ws.Range("F8").Value = 51.57091136640102
$ws.Range("G8").Value = 3.723581670176031
$ws.Range("L8").Value = 10.51187818208165
$ws.Range("M8").Value = 17.60030054677172
$ws.Range("B9").Value = 20.63015704556861
$ws.Range("C9").Value = 8.520993245698568
$ws.Range("D9").Value = 10.17925809528017
$ws.Range("F9").Value = 54.40460394223113
$ws.Range("G9").Value = 3.712274693924177
$ws.Range("L9").Value = 10.49844196832769
$ws.Range("M9").Value = 17.6457874711233
$ws.Range("B10").Value = 21.10014426433285
$ws.Range("C10").Value = 9.169186193110958
$ws.Range("D10").Value = 10.39855380416054
$ws.Range("F10").Value = 56.4292538686699
$ws.Range("G10").Value = 3.704671017659172
$ws.Range("L10").Value = 10.49492369243238
$ws.Range("M10").Value = 17.70310477664869
$ws.Range("B11").Value = 21.31989421058617
$ws.Range("C11").Value = 9.448727180373487
$ws.Range("D11").Value = 10.49706062845613
$ws.Range("F11").Value = 57.33447581630631
$ws.Range("G11").Value = 3.701362336983529
$ws.Range("L11").Value = 10.49470318562387
$ws.Range("M11").Value = 17.73433177721548
$ws.Range("B12").Value = 21.4038737753455
$ws.Range("C12").Value = 9.552378483493515
$ws.Range("D12").Value = 10.53416720122992
$ws.Range("F12").Value = 57.67472618076892
$ws.Range("G12").Value = 3.700130851863539
$ws.Range("L12").Value = 10.49481803799345
$ws.Range("M12").Value = 17.74689290048116
$ws.Range("B13").Value = 21.38575469456608
$ws.Range("C13").Value = 9.530153370086971
$ws.Range("D13").Value = 10.52618459357551
$ws.Range("F13").Value = 57.60156406398627
$ws.Range("G13").Value = 3.700395123327609
$ws.Range("L13").Value = 10.49478448259094
$ws.Range("M13").Value = 17.74415498527406
$ws.Range("B14").Value = 21.32678844781221
$ws.Range("C14").Value = 9.457298820022446
$ws.Range("D14").Value = 10.50011742221845
$ws.Range("F14").Value = 57.36252072913862
$ws.Range("G14").Value = 3.701260593227172
$ws.Range("L14").Value = 10.4947086599773
$ws.Range("M14").Value = 17.7353504662911
$ws.Range("B15").Value = 21.29076685281508
$ws.Range("C15").Value = 9.412386186359567
$ws.Range("D15").Value = 10.48412456284185
$ws.Range("F15").Value = 57.2157617603777
$ws.Range("G15").Value = 3.701793505713792
$ws.Range("L15").Value = 10.49468804446876
$ws.Range("M15").Value = 17.73005315057601
$ws.Range("B16").Value = 21.08589508882798
$ws.Range("C16").Value = 9.150609000054907
$ws.Range("D16").Value = 10.39208960051328
$ws.Range("F16").Value = 56.36975442279991
$ws.Range("G16").Value = 3.704890257163702
$ws.Range("L16").Value = 10.49496586564621
$ws.Range("M16").Value = 17.7011673164892
$ws.Range("B17").Value = 20.96167262937741
$ws.Range("C17").Value = 8.986090426329033
$ws.Range("D17").Value = 10.33529712787448
$ws.Range("F17").Value = 55.84651463153186
$ws.Range("G17").Value = 3.70682838197104
$ws.Range("L17").Value = 10.49548970112106
$ws.Range("M17").Value = 17.68476344678209
$ws.Range("B18").Value = 20.8907887193066
$ws.Range("C18").Value = 8.890021089746281
$ws.Range("D18").Value = 10.30251468846846
$ws.Range("F18").Value = 55.54408835297766
$ws.Range("G18").Value = 3.707957294272145
$ws.Range("L18").Value = 10.49592089572689
$ws.Range("M18").Value = 17.67581375696734
$ws.Range("B19").Value = 20.8668885344058
$ws.Range("C19").Value = 8.857245969547161
$ws.Range("D19").Value = 10.29139552600496
$ws.Range("F19").Value = 55.44144747596821
$ws.Range("G19").Value = 3.708341960690054
$ws.Range("L19").Value = 10.49608920266774
$ws.Range("M19").Value = 17.67286704378674
$ws.Range("B20").Value = 20.9748384438365
$ws.Range("C20").Value = 9.003753029162578
$ws.Range("D20").Value = 10.34135499911259
$ws.Range("F20").Value = 55.90236878797702
$ws.Range("G20").Value = 3.706620601591434
$ws.Range("L20").Value = 10.49542049440787
$ws.Range("M20").Value = 17.68645946386366
$ws.Range("B21").Value = 21.34408822495792
$ws.Range("C21").Value = 9.478757800969472
$ws.Range("D21").Value = 10.50777942089513
$ws.Range("F21").Value = 57.4328044426586
$ws.Range("G21").Value = 3.701005803191506
$ws.Range("L21").Value = 10.49472554860768
$ws.Range("M21").Value = 17.73791663178153
$ws.Range("B22").Value = 21.58983139639797
$ws.Range("C22").Value = 9.776350810733279
$ws.Range("D22").Value = 10.61539895455773
$ws.Range("F22").Value = 58.41813180140691
$ws.Range("G22").Value = 3.697461100968979
$ws.Range("L22").Value = 10.49542743960165
$ws.Range("M22").Value = 17.7758347256989
$ws.Range("B23").Value = 21.45829945857097
$ws.Range("C23").Value = 9.618695317334266
$ws.Range("D23").Value = 10.55807058195624
$ws.Range("F23").Value = 57.89369037300751
$ws.Range("G23").Value = 3.699341603145745
$ws.Range("L23").Value = 10.49494708806685
$ws.Range("M23").Value = 17.75520663089837
$ws.Range("B24").Value = 20.96888451133616
$ws.Range("C24").Value = 8.995772393056962
$ws.Range("D24").Value = 10.33861664450927
$ws.Range("F24").Value = 55.87712211442469
$ws.Range("G24").Value = 3.706714493430338
$ws.Range("L24").Value = 10.49545137772839
$ws.Range("M24").Value = 17.68569119562211
$ws.Range("B25").Value = 20.46295055549833
$ws.Range("C25").Value = 8.268723758565468
$ws.Range("D25").Value = 10.09749667874531
$ws.Range("F25").Value = 53.64685648519665
$ws.Range("G25").Value = 3.715209191032927
$ws.Range("L25").Value = 10.5009613730435
$ws.Range("M25").Value = 17.62927539008076
